# Changing the work report
# Adds the 29/7/2025(Onsite) 'Car Tracking Project' entries (rows 22-26) to Sheet1,
# reproducing the new shared-string insertion order from the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cell values, entered in the same order the strings were authored ---
# (this keeps the generated sharedStrings.xml insertion order stable)
$ws.Range("A22").Value = "29/7/2025(Onsite)"
$ws.Range("E22").Value = "Finding the right version of the code, although I added somehow good descriptive comment before each commit, need more descriptivness while writing `ndescriptions for each commit."
$ws.Range("C22").Value = "Recover the working version of the code."
$ws.Range("C24").Value = "Change the time of the autonomous extraction and test if it works on its own"
$ws.Range("C25").Value = "Let the data be saved automatically into the data subdirectory in the main after each run"
$ws.Range("C26").Value = "Start building the offline ai model that will clean the data in the data subdirectory then send it to a `nnew folder"
$ws.Range("C23").Value = "Test the pipeline locally"
$ws.Range("F22").Value = "DONE and the pipeline manual run on github is working perfectly الحمد الله`nThe run took around 18 to 19 minutes to finish"
$ws.Range("D23").Value = "The entry point in your code is run_pipeline.py from the main directory"
$ws.Range("F23").Value = "DONE: The run took around 18 to 21 minutes"

# --- Fill in the remaining Date / Project cells (reuse existing shared strings) ---
$projectVal = "Car Tracking Project"
$dateVal = "29/7/2025(Onsite)"
$ws.Range("B22").Value = $projectVal
$ws.Range("B23").Value = $projectVal
$ws.Range("B24").Value = $projectVal
$ws.Range("B25").Value = $projectVal
$ws.Range("B26").Value = $projectVal
$ws.Range("A23").Value = $dateVal
$ws.Range("A24").Value = $dateVal
$ws.Range("A25").Value = $dateVal
$ws.Range("A26").Value = $dateVal

# --- Wrap text for the cells that use the wrap style (matches style index 2) ---
$ws.Range("E22").WrapText = $true
$ws.Range("F22").WrapText = $true
$ws.Range("C26").WrapText = $true

# --- Row heights: rows 22 and 26 wrap to two lines (28.8pt = 2 x 14.4pt default) ---
$ws.Rows("22").RowHeight = 28.8
$ws.Rows("26").RowHeight = 28.8

# --- Final selection / active cell, matching the end of the editing session ---
$ws.Range("C26").Select()
